$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manager e-mail addresses (column C) to the new domain.
$ws.Range("C2").Value  = "gsoaresdesouzaaws+helena@gmail.com"
$ws.Range("C3").Value  = "gsoaresdesouzaaws+alice@gmail.com"
$ws.Range("C4").Value  = "gsoaresdesouzaaws+laura@gmail.com"
$ws.Range("C5").Value  = "gsoaresdesouzaaws+manuela@gmail.com"
$ws.Range("C6").Value  = "gsoaresdesouzaaws+valentina@gmail.com"
$ws.Range("C7").Value  = "gsoaresdesouzaaws+sophia@gmail.com"
$ws.Range("C8").Value  = "gsoaresdesouzaaws+isabella@gmail.com"
$ws.Range("C9").Value  = "gsoaresdesouzaaws+heloisa@gmail.com"
$ws.Range("C10").Value = "gsoaresdesouzaaws+luiza@gmail.com"
$ws.Range("C11").Value = "gsoaresdesouzaaws+julia@gmail.com"
$ws.Range("C12").Value = "gsoaresdesouzaaws+lorena@gmail.com"
$ws.Range("C13").Value = "gsoaresdesouzaaws+livia@gmail.com"
$ws.Range("C14").Value = "gsoaresdesouzaaws+maria_luiza@gmail.com"
$ws.Range("C15").Value = "gsoaresdesouzaaws+cecilia@gmail.com"
$ws.Range("C16").Value = "gsoaresdesouzaaws+eloa@gmail.com"
$ws.Range("C17").Value = "gsoaresdesouzaaws+miguel@gmail.com"
$ws.Range("C18").Value = "gsoaresdesouzaaws+arthur@gmail.com"
$ws.Range("C19").Value = "gsoaresdesouzaaws+heitor@gmail.com"
$ws.Range("C20").Value = "gsoaresdesouzaaws+bernardo@gmail.com"
$ws.Range("C21").Value = "gsoaresdesouzaaws+davi@gmail.com"
$ws.Range("C22").Value = "gsoaresdesouzaaws+theo@gmail.com"
$ws.Range("C23").Value = "gsoaresdesouzaaws+lorenzo@gmail.com"
$ws.Range("C24").Value = "gsoaresdesouzaaws+gabriel@gmail.com"
$ws.Range("C25").Value = "gsoaresdesouzaaws+pedro@gmail.com"
$ws.Range("C26").Value = "gsoaresdesouzaaws+benjamin@gmail.com"
$ws.Range("C27").Value = "gsoaresdesouzaaws+diretoria@gmail.com"

# The e-mails are no longer clickable mailto: hyperlinks.
$ws.Hyperlinks.Delete()

# Widen column C to fit the new (longer) addresses and drop the old "best fit" flag.
$ws.Columns.Item(3).ColumnWidth = 55.5

# Mark a helper cell (J23) with an underlined style, used while composing the e-mails.
$ws.Range("J23").Font.Underline = 2

# Printing options were configured for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the work finished.
$ws.Range("J23").Select()
